# Updates cryptos list values (Price / Volume(1h) columns) per the
# Fri Jun 23 18:12:26 UTC 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '31.191.35'
$ws.Range("E2").Value = '  +3.91%  '
$ws.Range("D3").Value = '1.921.48'
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5017'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3023'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06978'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.19%  '
$ws.Range("D10").Value = '1.914.04'
$ws.Range("E10").Value = '  +1.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07330'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '92.67'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6867'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.123'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.31%  '
$ws.Range("D16").Value = '31.142.00'
$ws.Range("E16").Value = '  +3.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008119'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.47'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +5.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("D20").Value = '2.160.59'
$ws.Range("E20").Value = '  +1.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9980'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.905'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.68%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '187.67'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +38.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.118'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.424'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.62'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +10.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.966'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.62%  '
$ws.Range("E29").Value = '  +1.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.381'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08997'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.087'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05271'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.66%  '
$ws.Range("E34").Value = '  +7.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.151'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.667'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01946'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +18.37%  '
$ws.Range("E38").Value = '  +2.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.215'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9441'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4403'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.965'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.51%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.916'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1345'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05886'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.674'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3913'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.407'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.24%  '
